# MPTrx - revised, has debug interface, and is ready for further testing
#
# Populates the new "Developer" / role-mapping debug columns (C:F) for the
# account rows (31-37) on the SIAM sheet, and leaves the selection on D32
# (the last cell edited), matching the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 34 role columns first (D/E), then the "Developer" marker filled down
# column C for every account row, then the remaining role in F34, and
# finally the "ExcelImporter" marker for the two accounts that have it.
$ws.Range("D34").Value2 = "AccountMgr"
$ws.Range("E34").Value2 = "PersonMgr"

$ws.Range("C31").Value2 = "Developer"
$ws.Range("C32").Value2 = "Developer"
$ws.Range("C33").Value2 = "Developer"
$ws.Range("C34").Value2 = "Developer"
$ws.Range("C35").Value2 = "Developer"
$ws.Range("C36").Value2 = "Developer"
$ws.Range("C37").Value2 = "Developer"

$ws.Range("F34").Value2 = "OrganizationMgr"

$ws.Range("D32").Value2 = "ExcelImporter"
$ws.Range("D33").Value2 = "ExcelImporter"

# Leave the selection where the author left it.
$ws.Range("D32").Select()
